$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New match rows (Open tournament results) ---
$rows = @(
    @{R=88; A=528; B=45816; D="Blois";              E=2;     F=3},
    @{R=89; A=329; B=45816; D="Blois";               E=$null; F=$null},
    @{R=90; A=421; B=45816; D="Blois Place"},
    @{R=91; A=557; B=45829; D="UCPA Paris 19";       E=2;     F=1},
    @{R=92; A=419; B=45829; D="UCPA Paris 19";       E=$null; F=$null},
    @{R=93; A=515; B=45829; D="UCPA Paris 19";       E=$null; F=$null},
    @{R=94; A=419; B=45836; D="National Squash 95";  E=2;     F=3},
    @{R=95; A=557; B=45836; D="National Squash 95";  E=$null; F=$null},
    @{R=96; A=516; B=45836; D="National Squash 95";  E=$null; F=$null},
    @{R=97; A=416; B=45843; D="Vincennes";           E=2;     F=2},
    @{R=98; A=416; B=45843; D="Vincennes Place";     E=$null; F=$null},
    @{R=99; A=476; B=45843; D="Vincennes"}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = "Open"
    $ws.Cells.Item($r, 4).Value = $row.D
}

$ws.Cells.Item(88, 5).Value = 2
$ws.Cells.Item(88, 6).Value = 3
$ws.Cells.Item(91, 5).Value = 2
$ws.Cells.Item(91, 6).Value = 1
$ws.Cells.Item(94, 5).Value = 2
$ws.Cells.Item(94, 6).Value = 3
$ws.Cells.Item(97, 5).Value = 2
$ws.Cells.Item(97, 6).Value = 2

# --- Replicate the existing row formatting (styles) down into the new rows ---
# Columns A:D formatting for every new row (88-99)
$ws.Range("A87:D87").Copy()
$ws.Range("A88:D99").PasteSpecial(-4122)  # xlPasteFormats

# Columns E:F formatting only for the rows that actually record a score
$ws.Range("E87:F87").Copy()
$ws.Range("E88:F89").PasteSpecial(-4122)
$ws.Range("E91:F98").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the view: selection on F99, scroll so row 85 is at the top ---
$ws.Range("F99").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 85
